$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1/J1 (style matches other headers - copy format from H1)
$ws.Range("H1").Copy()
$ws.Range("I1").Value = "I0"
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("J1").Value = "IF"
$ws.Range("J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Data for columns I (I0) and J (IF), rows 2-45
$data = @(
    @(7, 8),
    @(7, 7),
    @(7, 8),
    @(9, 9),
    @(4, 5),
    @(7, 7),
    @(7, 7),
    @(8, 8),
    @(8, 8),
    @(6, 7),
    @(8, 8),
    @(6, 6),
    @(9, 9),
    @(5, 6),
    @(1, 2),
    @(9, 9),
    @(8, 8),
    @(9, 9),
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(9, 9),
    @(8, 9),
    @(9, 9),
    @(10, 10),
    @(7, 8),
    @(3, 4),
    @(8, 9),
    @(7, 8),
    @(6, 7),
    @(7, 8),
    @(7, 7),
    @(4, 6),
    @(7, 8),
    @(8, 8),
    @(8, 8),
    @(5, 6),
    @(10, 10),
    @(6, 8),
    @(7, 8),
    @(4, 4),
    @(5, 6),
    @(7, 8),
    @(3, 3)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
